# Extend age tables for lookup completeness
#
# The "SoHPOtM" sheet holds a lookup table keyed by age (row 1 = age in
# years, row 2 = share available to open markets). It previously only
# covered ages 0-25 (columns B:AA); this extends it through age 210
# (columns AB:HD) so downstream lookups never fall off the end of the
# table. New header cells continue the existing alternating-alignment
# formatting (even ages right-aligned, odd ages default), and the new
# share cells continue the flat 0.5 plateau with the existing two-decimal
# number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoHPOtM")

$firstNewCol = 28   # column AB -> age 26
$lastCol     = 212  # column HD -> age 210

for ($col = $firstNewCol; $col -le $lastCol; $col++) {
    $age = $col - 2

    $headerCell = $ws.Cells.Item(1, $col)
    $headerCell.Value = $age
    if ($age % 2 -eq 0) {
        $headerCell.HorizontalAlignment = -4152   # xlRight, matches existing even-age columns
    }

    $shareCell = $ws.Cells.Item(2, $col)
    $shareCell.Value = 0.5
    $shareCell.NumberFormat = "0.00"
}

# Bring the lookup sheet to the front, with the newly-added corner of the
# table in view and selected, mirroring where the author ended up.
$ws.Activate()
$ws.Range("HD2").Select()
$excel.ActiveWindow.ScrollColumn = 188
$excel.ActiveWindow.ScrollRow = 1
